$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, pushing the existing rows 99-214 down
# to 100-215 (this also extends the sheet dimension to R215).
$ws.Rows(99).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A99").Value = 7
$ws.Range("B99").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C99").Value = "Ñuble"
$ws.Range("D99").Value = 44650
$ws.Range("E99").Value = 16
$ws.Range("F99").Value = 100112003
$ws.Range("G99").Value = "Ajo"
$ws.Range("H99").Value = "Chino"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 60
$ws.Range("K99").Value = 19000
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = 19500
$ws.Range("N99").Value = "$/caja 10 kilos"
$ws.Range("O99").Value = "China"
$ws.Range("P99").Value = 1950
$ws.Range("Q99").Value = 10
$ws.Range("R99").Value = "Hortaliza"
